$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remise en ordre des TD et TP du groupe 2 :
# on echange le contenu des lignes 97 et 98 (colonnes A, C, D, G, H)

# --- Colonnes texte simples : on permute directement les valeurs ---
$cols = @("A", "C", "D", "H")
foreach ($col in $cols) {
    $cellTop = $ws.Range("$col" + "97")
    $cellBottom = $ws.Range("$col" + "98")

    $valTop = $cellTop.Value2
    $valBottom = $cellBottom.Value2

    $cellTop.Value2 = $valBottom
    $cellBottom.Value2 = $valTop
}

# --- Colonne G : contient un nombre stocke en texte ("2") sur la ligne 98 ---
# et une cellule vide sur la ligne 97. On utilise Cut/Paste afin de
# preserver le type texte de la cellule deplacee.
$ws.Range("G98").Cut($ws.Range("G97"))
